# Actualización automática del tracker
# Fill in 'resultado' (G) / 'profit' (H) for rows whose bets have now
# settled, and fix A139/A140 (event_id) to be stored as numbers instead
# of text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (resultado, profit)
$updates = @{
    122 = @("Acierto", 1.1)
    130 = @("Fallo", -1)
    131 = @("Fallo", -1)
    132 = @("Fallo", -1)
    133 = @("Fallo", -1)
    134 = @("Acierto", 2.25)
    136 = @("Fallo", -1)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 7).Value = $vals[0]
    $ws.Cells.Item($row, 8).Value = $vals[1]
}

# A139 / A140 were stored as text; convert them to plain numbers.
$ws.Cells.Item(139, 1).Value = 14678166
$ws.Cells.Item(140, 1).Value = 14679464
